# Daily attendance processing - 2025-11-28 16:56:23
# Swap the order of the two comma-separated entries in the "Recorded By"
# column (G) wherever "dnasr281@gmail.com" is listed first, e.g.
#   "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
#   "dnasr281@gmail.com, admin@admin.com"  -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val) {
        $parts = $val -split ', '
        if ($parts.Count -eq 2 -and $parts[0] -eq 'dnasr281@gmail.com') {
            $cell.Value = "$($parts[1]), $($parts[0])"
        }
    }
}
